# Update "想去人数" (interest count) figures that changed between crawls.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1108
$wsExpo.Range("F4").Value = 1768
$wsExpo.Range("F5").Value = 786
$wsExpo.Range("F6").Value = 240
$wsExpo.Range("F7").Value = 205

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 13

# Sheet "全部类型" (All types, combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1108
$wsAll.Range("F4").Value = 1768
$wsAll.Range("F5").Value = 13
$wsAll.Range("F6").Value = 786
$wsAll.Range("F7").Value = 240
$wsAll.Range("F8").Value = 205
